$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current layout:
#   B24/C24 -> "LOM3254 -  Laboratório de Circuitos Elétricos  (Indicação de Conjunto)\n"
#   B25/C25 -> "LOB1053 -  Física III  (Requisito)\n"
# Target layout (swap the two requirement lines):
#   B24/C24 -> "LOB1053 -  Física III  (Requisito)\n"
#   B25/C25 -> "LOM3254 -  Laboratório de Circuitos Elétricos  (Indicação de Conjunto)\n"

$reqLOB1053 = "LOB1053 -  Física III  (Requisito)`n"
$reqLOM3254 = "LOM3254 -  Laboratório de Circuitos Elétricos  (Indicação de Conjunto)`n"

$ws.Range("B24").Value = $reqLOB1053
$ws.Range("C24").Value = $reqLOB1053

$ws.Range("B25").Value = $reqLOM3254
$ws.Range("C25").Value = $reqLOM3254
